$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1007, shifting existing rows 1007-1093 down to 1008-1094
$ws.Rows("1007:1007").Insert()

# Populate the newly inserted row with the new data record
$ws.Range("A1007").Value = 6
$ws.Range("B1007").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1007").Value = "Metropolitana"
$ws.Range("D1007").Value = 44918
$ws.Range("E1007").Value = 13
$ws.Range("F1007").Value = 100112021
$ws.Range("G1007").Value = "Ají"
$ws.Range("H1007").Value = "Americana (o)"
$ws.Range("I1007").Value = "Primera"
$ws.Range("J1007").Value = 470
$ws.Range("K1007").Value = 20000
$ws.Range("L1007").Value = 27000
$ws.Range("M1007").Value = 23191
$ws.Range("N1007").Value = "`$/caja 15 kilos"
$ws.Range("O1007").Value = "Región de O'Higgins"
$ws.Range("P1007").Value = 1546
$ws.Range("Q1007").Value = 15
$ws.Range("R1007").Value = "Hortaliza"
